$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Extent of Contamination")

# Row 4: Indoor / Area Contaminated -> Area value reset to 0
$ws.Range("G4").Value = 0

# Row 5: Indoor / Loading -> distribution type switched from Uniform to Constant
$ws.Range("F5").Value = "Constant"
$ws.Range("G5").Value = -1
$ws.Range("H5").ClearContents()

# Row 6: Underground / Area Contaminated -> Area value set
$ws.Range("G6").Value = 2682.8539999999998

# Row 7: Underground / Loading -> distribution type switched from Constant to Uniform
$ws.Range("F7").Value = "Uniform"
$ws.Range("G7").Value = 5.6318533337268804
$ws.Range("H7").Value = 6.1323408067707499
$ws.Range("G7").ClearFormats()

# Rows 8-10: Indoor Contamination Breakout (Residential/Commercial/Industrial) -> zeroed out
$ws.Range("G8").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("G10").Value = 0

# Rows 21-26: Underground Surface Type Breakout -> new fractions
$ws.Range("G21").Value = 0.5
$ws.Range("G22").Value = 0.125
$ws.Range("G24").Value = 0.125
$ws.Range("G25").Value = 0.25

# Rows 27-32: Indoor Surface Type Breakout -> zeroed out
$ws.Range("G27").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("G30").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("G32").Value = 0

# Update the active selection/view on this sheet
$ws.Range("E9").Select()
